# Clean up the "Tipo" column (column I) on the single worksheet:
# many near-duplicate / inconsistent category labels (and a few cached
# #N/A errors) are consolidated into the canonical labels already used
# elsewhere in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = "GRANEL"
    13  = "GRANEL"
    16  = "GRANEL"
    23  = "CINTA"
    28  = "EXPORTAÇÃO"
    32  = "EXPORTAÇÃO"
    34  = "GRANEL"
    35  = "GRANEL"
    44  = "GRANEL"
    46  = "GRANEL"
    55  = "GRANEL"
    56  = "GRANEL"
    62  = "BALDE "
    63  = "BALDE "
    65  = "GRANEL"
    68  = "GRANEL"
    73  = "GRANEL"
    74  = "EXPORTAÇÃO"
    76  = "GRANEL"
    77  = "GRANEL"
    84  = "GRANEL"
    89  = "GRANEL"
    90  = "GRANEL"
    91  = "CINTA"
    93  = "EXPORTAÇÃO"
    94  = "EXPORTAÇÃO"
    95  = "EXPORTAÇÃO"
    96  = "INDUSTRIA"
    133 = "LÍQUIDO"
    134 = "LÍQUIDO"
    135 = "LÍQUIDO"
    138 = "KIT "
    139 = "KIT "
    140 = "KIT "
    141 = "CINTA"
    142 = "KIT"
    143 = "CINTA"
    144 = "EXPORTAÇÃO"
    146 = "KIT "
    149 = "KIT "
    152 = "EMBALADO"
    161 = "GRANEL"
    164 = "EXPORTAÇÃO"
    173 = "CINTA"
    179 = "CINTA"
    182 = "OVO EM PÓ"
    183 = "OVO EM PÓ"
    184 = "LÍQUIDO"
    185 = "LÍQUIDO"
    186 = "GRANEL"
    190 = "KIT"
    191 = "KIT"
    192 = "KIT"
    193 = "KIT"
    194 = "KIT"
    195 = "KIT"
    196 = "KIT"
    198 = "GRANEL"
    199 = "GRANEL"
    200 = "KIT"
    201 = "Industria"
    211 = "GRANEL"
    217 = "EMBALADO"
    218 = "GRANEL"
    219 = "GRANEL"
    230 = "EXPORTAÇÃO"
    233 = "EXPORTAÇÃO"
    234 = "EXPORTAÇÃO"
    235 = "EXPORTAÇÃO"
    236 = "EXPORTAÇÃO"
    237 = "EXPORTAÇÃO"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 9).Value = $updates[$row]
}

$ws.Range("I187").Select()
